$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect so cell values can be edited
$ws.Unprotect()

# Update the confidential disclosure text date: 2021-07-07 -> 2021-07-08
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-08 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-10
$ws.Range("D2").Value = 0.1358856333797115
$ws.Range("E2").Value = -0.0002254791431790526

$ws.Range("D3").Value = 0.1082675523825937
$ws.Range("E3").Value = -0.01172832302485693

$ws.Range("D4").Value = 0.1131698525721662
$ws.Range("E4").Value = -0.00821005917159745

$ws.Range("D5").Value = 0.1190462954025016
$ws.Range("E5").Value = -0.01069559864573721

$ws.Range("D6").Value = 0.1219396632978936
$ws.Range("E6").Value = -0.006536819637139701

$ws.Range("D7").Value = 0.1428445622119583
$ws.Range("E7").Value = -0.01096575606002292

$ws.Range("D8").Value = 0.1306819933460119
$ws.Range("E8").Value = -0.0117252931323284

$ws.Range("D9").Value = 0.1281644474071633
$ws.Range("E9").Value = -0.007148661197956985

$ws.Range("E10").Value = -0.008314823886128919

# Restore sheet protection so the worksheet remains protected, as it was originally.
$ws.Protect()

